$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows above row 285 (new data week), shifting the rest of the
# table down by 4 rows (old row 285 -> new row 289, ..., old row 328 -> new row 332).
$ws.Range("A285:A288").EntireRow.Insert()

# Constant columns shared by every data row in this sheet.
$mercadoId = 6
$mercado = "Mercado Mayorista Lo Valledor de Santiago"
$region = "Metropolitana"
$codreg = 13
$categoriaId = 100112043
$categoria = "Pepino dulce"
$variedad = "Cultivar IV Región"
$unidad = "`$/bandeja 18 kilos"
$origen = "Provincia de Limarí"
$kgUnidades = 18
$clasificacion = "Hortaliza"

# New week of data (Fecha serial 45015) with its four quality rows.
$rows = @(
    @{ Row = 285; Calidad = "Especial"; Volumen = 320; PMin = 14000; PMax = 14000; PProm = 14000; PKg = 778 },
    @{ Row = 286; Calidad = "Primera";  Volumen = 910; PMin = 11000; PMax = 12000; PProm = 11626; PKg = 646 },
    @{ Row = 287; Calidad = "Segunda";  Volumen = 630; PMin = 9000;  PMax = 10000; PProm = 9429;  PKg = 524 },
    @{ Row = 288; Calidad = "Tercera";  Volumen = 230; PMin = 7000;  PMax = 7000;  PProm = 7000;  PKg = 389 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $mercadoId
    $ws.Cells.Item($row, 2).Value = $mercado
    $ws.Cells.Item($row, 3).Value = $region
    $ws.Cells.Item($row, 4).Value = 45015
    $ws.Cells.Item($row, 5).Value = $codreg
    $ws.Cells.Item($row, 6).Value = $categoriaId
    $ws.Cells.Item($row, 7).Value = $categoria
    $ws.Cells.Item($row, 8).Value = $variedad
    $ws.Cells.Item($row, 9).Value = $r.Calidad
    $ws.Cells.Item($row, 10).Value = $r.Volumen
    $ws.Cells.Item($row, 11).Value = $r.PMin
    $ws.Cells.Item($row, 12).Value = $r.PMax
    $ws.Cells.Item($row, 13).Value = $r.PProm
    $ws.Cells.Item($row, 14).Value = $unidad
    $ws.Cells.Item($row, 15).Value = $origen
    $ws.Cells.Item($row, 16).Value = $r.PKg
    $ws.Cells.Item($row, 17).Value = $kgUnidades
    $ws.Cells.Item($row, 18).Value = $clasificacion
}
